$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 211, shifting the existing rows 211-219 down to 214-222.
$ws.Range("A211:T213").EntireRow.Insert()

# Populate the 3 newly inserted rows (211-213) with a new weekly price block.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant for this product/vendor and are
# copied from the surrounding rows; D,L,M,N,O,P,S carry the new data.

# Row 211 - Especial
$ws.Range("A211").Value = 2
$ws.Range("B211").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44448
$ws.Range("E211").Value = 4
$ws.Range("F211").Value = "Fruta"
$ws.Range("G211").Value = 100101
$ws.Range("H211").Value = "Berries"
$ws.Range("I211").Value = 100112025
$ws.Range("J211").Value = "Frutilla"
$ws.Range("K211").Value = "Sin especificar"
$ws.Range("L211").Value = "Especial"
$ws.Range("M211").Value = 200
$ws.Range("N211").Value = 30000
$ws.Range("O211").Value = 31000
$ws.Range("P211").Value = 30500
$ws.Range("Q211").Value = "$/bandeja 7 kilos"
$ws.Range("R211").Value = "Provincia de Melipilla"
$ws.Range("S211").Value = 4357
$ws.Range("T211").Value = 7

# Row 212 - Primera
$ws.Range("A212").Value = 2
$ws.Range("B212").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C212").Value = "Coquimbo"
$ws.Range("D212").Value = 44448
$ws.Range("E212").Value = 4
$ws.Range("F212").Value = "Fruta"
$ws.Range("G212").Value = 100101
$ws.Range("H212").Value = "Berries"
$ws.Range("I212").Value = 100112025
$ws.Range("J212").Value = "Frutilla"
$ws.Range("K212").Value = "Sin especificar"
$ws.Range("L212").Value = "Primera"
$ws.Range("M212").Value = 300
$ws.Range("N212").Value = 25000
$ws.Range("O212").Value = 26000
$ws.Range("P212").Value = 25500
$ws.Range("Q212").Value = "$/bandeja 7 kilos"
$ws.Range("R212").Value = "Provincia de Melipilla"
$ws.Range("S212").Value = 3643
$ws.Range("T212").Value = 7

# Row 213 - Segunda
$ws.Range("A213").Value = 2
$ws.Range("B213").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C213").Value = "Coquimbo"
$ws.Range("D213").Value = 44448
$ws.Range("E213").Value = 4
$ws.Range("F213").Value = "Fruta"
$ws.Range("G213").Value = 100101
$ws.Range("H213").Value = "Berries"
$ws.Range("I213").Value = 100112025
$ws.Range("J213").Value = "Frutilla"
$ws.Range("K213").Value = "Sin especificar"
$ws.Range("L213").Value = "Segunda"
$ws.Range("M213").Value = 240
$ws.Range("N213").Value = 20000
$ws.Range("O213").Value = 21000
$ws.Range("P213").Value = 20500
$ws.Range("Q213").Value = "$/bandeja 7 kilos"
$ws.Range("R213").Value = "Provincia de Melipilla"
$ws.Range("S213").Value = 2929
$ws.Range("T213").Value = 7
